# Weekly update: a new observation is inserted at row 299 of the
# "Berenjena" (Mercado Mayorista Lo Valledor de Santiago) sheet, pushing
# the previously-existing rows 299-346 down to 300-347.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 299 (shifts 299..346 -> 300..347)
$ws.Rows(299).Insert()

# Populate the newly inserted row 299 with the new weekly record.
$ws.Cells.Item(299, 1).Value = 6
$ws.Cells.Item(299, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(299, 3).Value = "Metropolitana"
$ws.Cells.Item(299, 4).Value = 45154
$ws.Cells.Item(299, 5).Value = 13
$ws.Cells.Item(299, 6).Value = 100112001
$ws.Cells.Item(299, 7).Value = "Berenjena"
$ws.Cells.Item(299, 8).Value = "Sin especificar"
$ws.Cells.Item(299, 9).Value = "Primera"
$ws.Cells.Item(299, 10).Value = 220
$ws.Cells.Item(299, 11).Value = 7000
$ws.Cells.Item(299, 12).Value = 8000
$ws.Cells.Item(299, 13).Value = 7455
$ws.Cells.Item(299, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(299, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(299, 16).Value = 149
$ws.Cells.Item(299, 17).Value = 50
$ws.Cells.Item(299, 18).Value = "Hortaliza"
